# spent_data_master_2_2020.xlsx - refactoring edit
#
# - A33 label: "Pre 19-20 CDEL Forecast Total" -> "Pre-profile CDEL Forecast one off new costs"
# - A36 label: "Pre 19-20 Forecast Non-Gov"     -> "Pre-profile Forecast Non-Gov"
# - C33 value: 0 -> 200
# - A33 picks up an explicit (re-applied) cell style
# - row 33 height / column A width tweak, selection moves to C39

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- label text updates (shared-string content) ---
$ws.Range("A33").Value = "Pre-profile CDEL Forecast one off new costs"
$ws.Range("A36").Value = "Pre-profile Forecast Non-Gov"

# --- updated forecast figure ---
$ws.Range("C33").Value = 200

# --- A33 now carries its own explicit style (re-applying the Normal style
#     stamps a dedicated cellXfs entry instead of the default/shared one) ---
$ws.Range("A33").Style = "Normal"

# --- row/column sizing tweaks ---
$ws.Rows.Item(33).RowHeight = 13.8
$ws.Columns.Item(1).ColumnWidth = 41.92

# --- move the active selection / scroll position ---
$ws.Range("C39").Select()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
